$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a NATMI ligand-receptor pair export (Fgf21 -> Fgfr4).
# The update re-runs the analysis with a refreshed TPM table: the existing
# 4 rows (FAPs -> {ECs, FAPs, MuSCs, Resolving-Mac}) get new numeric values,
# the old "MuSCs" target becomes "Inflammatory-Mac", a new "MuSCs" target
# row is added for the FAPs sender, and a whole new block of 6 rows is
# appended for a second sender cluster "MuSCs" -> Fgfr4.
# ---------------------------------------------------------------------------

$rows = @(
    @("FAPs","Fgf21","Fgfr4","ECs",1,0.3333333333333333,0.08847333333333333,0.26542,0.4424168072381536,0.5434168393294406,3,1,0.363908,1.091724,0.01118972054940699,0.01663265480083782,0.03219615378666666,0.28976538408,0.0049505204393558,0.009038464701528935),
    @("FAPs","Fgf21","Fgfr4","FAPs",1,0.3333333333333333,0.08847333333333333,0.26542,0.4424168072381536,0.5434168393294406,1,0.3333333333333333,0.153566,0.460698,0.004721964413781051,0.007018835164781924,0.01358649590666667,0.12227846316,0.002089076419837192,0.003814153221020127),
    @("FAPs","Fgf21","Fgfr4","Inflammatory-Mac",1,0.3333333333333333,0.08847333333333333,0.26542,0.4424168072381536,0.5434168393294406,1,0.3333333333333333,0.06252866666666666,0.187586,0.001922679101110775,0.002857913889838424,0.005532119568888889,0.04978907612,0.0008506255492569525,0.001553038533091703),
    @("FAPs","Fgf21","Fgfr4","MuSCs",1,0.3333333333333333,0.08847333333333333,0.26542,0.4424168072381536,0.5434168393294406,2,1,31.927516,63.85503199999999,0.9817315966582778,0.9728454303033116,2.824733765573333,16.94840259344,0.4343345585583701,0.5286605888915152),
    @("FAPs","Fgf21","Fgfr4","Resolving-Mac",1,0.3333333333333333,0.08847333333333333,0.26542,0.4424168072381536,0.5434168393294406,1,0.3333333333333333,0.01411566666666667,0.042347,0.0004340392774233579,0.0006451658412300904,0.001248860082222222,0.01123974074,0.0001920262713335972,0.0003505939822845755),
    @("MuSCs","Fgf21","Fgfr4","ECs",1,0.5,0.111504,0.223008,0.5575831927618463,0.4565831606705595,3,1,0.363908,1.091724,0.01118972054940699,0.01663265480083782,0.040577197632,0.243463185792,0.006239200110051192,0.007594190099308887),
    @("MuSCs","Fgf21","Fgfr4","FAPs",1,0.5,0.111504,0.223008,0.5575831927618463,0.4565831606705595,1,0.3333333333333333,0.153566,0.460698,0.004721964413781051,0.007018835164781924,0.017123223264,0.102739339584,0.002632887993943858,0.003204681943761798),
    @("MuSCs","Fgf21","Fgfr4","Inflammatory-Mac",1,0.5,0.111504,0.223008,0.5575831927618463,0.4565831606705595,1,0.3333333333333333,0.06252866666666666,0.187586,0.001922679101110775,0.002857913889838424,0.006972196448,0.041833178688,0.001072053551853823,0.001304875356746721),
    @("MuSCs","Fgf21","Fgfr4","MuSCs",1,0.5,0.111504,0.223008,0.5575831927618463,0.4565831606705595,2,1,31.927516,63.85503199999999,0.9817315966582778,0.9728454303033116,3.560045744064,14.240182976256,0.5473970380999077,0.4441848414117965),
    @("MuSCs","Fgf21","Fgfr4","Resolving-Mac",1,0.5,0.111504,0.223008,0.5575831927618463,0.4565831606705595,1,0.3333333333333333,0.01411566666666667,0.042347,0.0004340392774233579,0.0006451658412300904,0.001573953296,0.009443719776000002,0.0002420130060897607,0.0002945718589455151)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
